$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)
$rows = @(
    @{ Row = 2;  D = 44455; M = 200; N = 12000; O = 13000; P = 12500; S = 6250 },
    @{ Row = 3;  D = 44461; M = 200; N = 11000; O = 12000; P = 11500; S = 5750 },
    @{ Row = 4;  D = 44489; M = 160; N = 9500;  O = 10000; P = 9750;  S = 4875 },
    @{ Row = 5;  D = 44482; M = 240; N = 10000; O = 11000; P = 10500; S = 5250 },
    @{ Row = 6;  D = 44454; M = 160; N = 12000; O = 13000; P = 12500; S = 6250 },
    @{ Row = 7;  D = 44490; M = 400; N = 9500;  O = 10000; P = 9750;  S = 4875 },
    @{ Row = 8;  D = 44497; M = 500; N = 9000;  O = 10000; P = 9500;  S = 4750 },
    @{ Row = 9;  D = 44475; M = 240; N = 11000; O = 12000; P = 11500; S = 5750 },
    @{ Row = 10; D = 44517; M = 400; N = 5500;  O = 6000;  P = 5750;  S = 2875 }
)

foreach ($r in $rows) {
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("M" + $r.Row).Value = $r.M
    $ws.Range("N" + $r.Row).Value = $r.N
    $ws.Range("O" + $r.Row).Value = $r.O
    $ws.Range("P" + $r.Row).Value = $r.P
    $ws.Range("S" + $r.Row).Value = $r.S
}
